# Workbook re-scraped on 27-11-2023 20:21: a handful of fixtures had
# their home/away rows reordered within their shared kickoff slot, and
# one newly-played match was appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    # Positional params only -- named (-RowA/-RowB) binding isn't
    # reliable in this host's PowerShell subset.
    param([int]$RowA, [int]$RowB)

    # Only columns F:V (home..url_partida) are swapped; A:E (Indice/
    # pais/torneio/temporada/data_partida) stay put since both fixtures
    # share that same row slot/date.
    $rngA = "F" + $RowA + ":V" + $RowA
    $rngB = "F" + $RowB + ":V" + $RowB

    $valsA = $ws.Range($rngA).Value()
    $valsB = $ws.Range($rngB).Value()

    $ws.Range($rngB).Value = $valsA
    $ws.Range($rngA).Value = $valsB
}

Swap-MatchRows 6 7
Swap-MatchRows 21 22
Swap-MatchRows 44 45

# Append the new fixture as row 118 (Indice 117). Copy formats from the
# previously-last data row (117) for the styled columns (A: bold+border
# index cell, E: datetime-formatted date cell) so no new style entries
# are created, then fill in the values.
$ws.Range("A117").Copy()
$ws.Range("A118").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E117").Copy()
$ws.Range("E118").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A118").Value = 117
$ws.Range("B118").Value = "turkey"
$ws.Range("C118").Value = "1-lig"
$ws.Range("D118").Value = "2023-2024"
$ws.Range("E118").Value = 45257.75
$ws.Range("F118").Value = "Genclerbirligi"
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = "Keciorengucu"
$ws.Range("I118").Value = 1
$ws.Range("J118").Value = 1.88
$ws.Range("K118").Value = "20/11/2023 16:12"
$ws.Range("L118").Value = 2.19
$ws.Range("M118").Value = "27/11/2023 17:57"
$ws.Range("N118").Value = 3.55
$ws.Range("O118").Value = "20/11/2023 16:12"
$ws.Range("P118").Value = 3.36
$ws.Range("Q118").Value = "27/11/2023 17:57"
$ws.Range("R118").Value = 4.13
$ws.Range("S118").Value = "20/11/2023 16:12"
$ws.Range("T118").Value = 3.44
$ws.Range("U118").Value = "27/11/2023 17:57"
$ws.Range("V118").Value = "https://www.betexplorer.com/football/turkey/1-lig/genclerbirligi-keciorengucu/OhYZl0vq/"
